# Replace curly double quotation marks (“ ”) surrounding quoted speech in the
# English (column C) dialogue lines with straight single quotes ('...').
# NOTE: for lines whose new text starts with an apostrophe, Excel treats a
# leading ' as a "store as text" prefix character and strips it from the
# stored value, so we double it (''...) to make one survive in the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(36, 3).Value = "[name=""Recording""]  'Chief Wei, you must act immediately upon receiving this message.'`n"
$ws.Cells.Item(38, 3).Value = "[name=""Recording""]  'The following has been Arts encrypted.'`n"
$ws.Cells.Item(39, 3).Value = "[name=""Deep Male Voice""]  'Mr. Wei, Groups Three and Four are standing down. We are unable to get to the source of the matter.'`n"
$ws.Cells.Item(40, 3).Value = "[name=""Deep Male Voice""]  'I can’t even speculate on who the mastermind might be.'`n"
$ws.Cells.Item(41, 3).Value = "[name=""Deep Male Voice""]  'They were whispering from the parliamentary gallery, laughing at my miserable, helpless performance. But I found no evidence to hold them accountable.'`n"
$ws.Cells.Item(42, 3).Value = "[name=""Deep Male Voice""]  'If something happens in Chernobog, you must do everything you can to stop it. Otherwise, who knows what will happen?'`n"
$ws.Cells.Item(43, 3).Value = "[name=""Deep Male Voice""]  'You have the wisdom and ability to stop it at its source. What we could not do, you must do in our stead.'`n"
$ws.Cells.Item(48, 3).Value = "[name=""Deep Male Voice""]  'I was unable to meet Chairman Witte, and the Messenger he sent to liaise with me was waylaid by unknown forces. Fortunately, that Messenger is still safe and sound.'`n"
$ws.Cells.Item(49, 3).Value = "[name=""Deep Male Voice""]  'The Messenger slipped out of Deity Grypherburg in the dead of night, and someone secretly aided him. I believe there are forces within Ursus vying for supremacy.'`n"
$ws.Cells.Item(50, 3).Value = "[name=""Deep Male Voice""]  'I was attacked many times on the road after that. And there were many others who tried to protect me.'`n"
$ws.Cells.Item(52, 3).Value = "[name=""Deep Male Voice""]  'I reached the Ural Rift and commandeered a messaging station. What happens now, even I don’t know.'`n"
$ws.Cells.Item(53, 3).Value = "[name=""Deep Male Voice""]  'I hope to make it back to Lungmen, to enjoy the flavors of home once again.'`n"
$ws.Cells.Item(117, 3).Value = "[name=""Wei Yenwu""]  'What did they do?'`n"
$ws.Cells.Item(155, 3).Value = "[name=""Ch’en""]  When you say 'act,' do you mean——`n"
$ws.Cells.Item(183, 3).Value = "''Don’t feel bad. I understand. We’re sworn brothers, aren’t we? Brothers... remember where we came from.'`n"
$ws.Cells.Item(186, 3).Value = "''I hate you. And I hate them. I should love you, but right now I hate every last one of you.'`n"
$ws.Cells.Item(189, 3).Value = "''Why me? Why does it have to be me? Who could feel safe sitting here?'`n"
$ws.Cells.Item(196, 3).Value = "[name=""Ch’en""]  The only difference is, I’m not like you. I don’t think anyone is a 'mistake.'`n"
$ws.Cells.Item(241, 3).Value = "''Brother, if this sword is used to slay dragons, does it also work on Dracos?'`n"
$ws.Cells.Item(242, 3).Value = "''It just might. Looks like I’ll have to be more careful around you, haha...'`n"
